$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark around "Ideal "/"Requirements" at the top.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2) Insert two new paragraphs right after the "To ensure APIs..." paragraph
#    and before the "ToDos" heading: an empty bordered paragraph, followed
#    by a bordered paragraph with the new note text (with a mid-run
#    "_GoBack" bookmark and a proofed "apis" word).
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.Text = "To ensure "
$f.Execute() | Out-Null
$r0 = $d.Range($f.Parent.Start, $f.Parent.Start)
$toEnsurePara = $r0.Paragraphs(1)

# Append an empty paragraph right after "To ensure ..." and turn it into
# the bordered blank paragraph.
$toEnsurePara.Range.InsertParagraphAfter()
$newEmptyPara = $toEnsurePara.Next()

$r1 = $newEmptyPara.Range
$r1.Collapse(1)
$emptyBorderedXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
"@
$r1.InsertXML($emptyBorderedXml)

# Append a second new paragraph right after the blank bordered one and fill
# it in with the note text.
$blankPara = $toEnsurePara.Next()
$blankPara.Range.InsertParagraphAfter()
$introPara = $blankPara.Next()

$r2 = $introPara.Range
$r2.Collapse(1)
$introXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Introducing external </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>apis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> adds more points </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>of failure to the system</w:t></w:r></w:p>
"@
$r2.InsertXML($introXml)

# ---------------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the "Find out how to not hard
#    code urls" run onto the "ToDos" run.
# ---------------------------------------------------------------------------
$todosFind = $d.Content.Find
$todosFind.Text = "ToDos"
$todosFind.Execute() | Out-Null
$todosR0 = $d.Range($todosFind.Parent.Start, $todosFind.Parent.Start)
$todosPara = $todosR0.Paragraphs(1)

$todosRange = $todosPara.Range
$todosRange.Collapse(1)
$todosXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="14B99B23" w14:textId="1A4294B2" w:rsidR="00514CDA" w:rsidRPr="00514CDA" w:rsidRDefault="00514CDA" w:rsidP="00C45C2E"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00514CDA"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>ToDos</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
$todosRange.InsertXML($todosXml)

$urlsFind = $d.Content.Find
$urlsFind.Text = "Find out how to not hard code "
$urlsFind.Execute() | Out-Null
$urlsR0 = $d.Range($urlsFind.Parent.Start, $urlsFind.Parent.Start)
$urlsPara = $urlsR0.Paragraphs(1)

$urlsRange = $urlsPara.Range
$urlsRange.Collapse(1)
$urlsXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="72766724" w14:textId="13313049" w:rsidR="00D74D31" w:rsidRDefault="00D74D31" w:rsidP="003E5138"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Find out how to not hard code </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>urls</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
$urlsRange.InsertXML($urlsXml)

Write-Output "done"
